# Regenerate the "K" (strikeouts) column (column G) values in save_data
# sheet, replacing the previous "Strike#" derived values with the
# recalculated K values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 0
    3  = 2
    4  = 1
    5  = 2
    6  = 3
    7  = 1
    8  = 7
    9  = 2
    10 = 2
    11 = 0
    12 = 0
    13 = 1
    14 = 2
    15 = 1
    16 = 1
    17 = 2
    18 = 1
    19 = 0
    20 = 1
    21 = 2
    22 = 1
    23 = 2
    24 = 0
    25 = 2
    26 = 0
    27 = 3
    28 = 2
    29 = 3
    30 = 2
    31 = 4
    32 = 1
    33 = 1
    34 = 1
    35 = 1
    36 = 5
    37 = 2
    38 = 0
    39 = 1
    40 = 3
    41 = 1
    42 = 1
    43 = 2
    45 = 1
    46 = 2
    47 = 1
    48 = 0
    49 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
